$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.509.64'
$ws.Range("E2").Value = '  +0.33%  '

# Row 3
$ws.Range("D3").Value = '1.847.66'
$ws.Range("E3").Value = '  +0.20%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.31%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '262.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.78%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.13%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5213'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.79%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3241'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.55%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06759'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.48%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7708'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.20%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07746'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '

# Row 13
$ws.Range("D13").Value = '1.860.01'
$ws.Range("E13").Value = '  +0.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.006'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.46%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.40%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.41%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007920'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("D20").Value = '26.569.17'
$ws.Range("E20").Value = '  +0.40%  '

# Row 21
$ws.Range("D21").Value = '2.091.90'
$ws.Range("E21").Value = '  -0.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.618'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.72%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.462'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.970'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.30%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.184'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.688'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.74%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.166'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.92%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08747'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.40%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.103'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.83%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04802'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.59%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.877'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.13%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.63%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7115'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.10%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.101'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01782'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.70%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.182'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4833'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.96%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.10%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8965'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.054'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.84%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.19%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.627'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05908'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.13%  '

# Row 47
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4139'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.02%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.040'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.08%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1229'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.31%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.05%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8819'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.42%  '
